$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.545.00'
$ws.Range("E2").Value = '  +1.31%  '
$ws.Range("D3").Value = '1.913.21'
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.59'
$ws.Range("E5").Value = '  +4.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.662'
$ws.Range("E6").Value = '  +6.22%  '
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.30'
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.348'
$ws.Range("E9").Value = '  +5.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '49.32'
$ws.Range("E10").Value = '  +5.74%  '
$ws.Range("E11").Value = '  +3.21%  '
$ws.Range("E12").Value = '  +0.99%  '
$ws.Range("E13").Value = '  +2.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.40'
$ws.Range("E14").Value = '  +8.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.704'
$ws.Range("E15").Value = '  +3.67%  '
$ws.Range("D16").Value = '1.912.81'
$ws.Range("E16").Value = '  +2.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.88'
$ws.Range("E17").Value = '  +3.80%  '
$ws.Range("D18").Value = '35.556.24'
$ws.Range("E18").Value = '  +1.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.64'
$ws.Range("E19").Value = '  +3.33%  '
$ws.Range("E20").Value = '  +4.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '244.95'
$ws.Range("E22").Value = '  +4.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.86'
$ws.Range("E23").Value = '  +2.34%  '
$ws.Range("E24").Value = '  +0.27%  '
$ws.Range("E25").Value = '  +1.36%  '
$ws.Range("E26").Value = '  +17.49%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '171.76'
$ws.Range("E27").Value = '  +0.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.50'
$ws.Range("E28").Value = '  +6.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.35'
$ws.Range("E29").Value = '  +3.68%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.130'
$ws.Range("E30").Value = '  +4.13%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.980'
$ws.Range("E31").Value = '  +25.30%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.17'
$ws.Range("E32").Value = '  +4.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0573'
$ws.Range("E33").Value = '  +2.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.23'
$ws.Range("E34").Value = '  +4.95%  '
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.73'
$ws.Range("E36").Value = '  +6.44%  '
$ws.Range("E37").Value = '  +1.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.34'
$ws.Range("E38").Value = '  +3.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.12'
$ws.Range("E39").Value = '  +3.28%  '
$ws.Range("E40").Value = '  +1.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '92.66'
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0638'
$ws.Range("E42").Value = '  +16.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '15.66'
$ws.Range("E43").Value = '  +5.26%  '
$ws.Range("D44").Value = '1.352.40'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("E45").Value = '  +2.60%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '47.64'
$ws.Range("E46").Value = '  +39.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.65'
$ws.Range("E47").Value = '  -1.27%  '
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.80'
$ws.Range("E49").Value = '  +2.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.60'
$ws.Range("E50").Value = '  +2.48%  '
$ws.Range("D51").Value = '2.097.62'
$ws.Range("E51").Value = '  +2.91%  '
